# Updates the logistic-regression odds-ratio table (columns A:C, rows 2-31)
# to reflect the retrained/re-selected model's coefficients and ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Index (col A), Feature (col B), Odds Ratio (col C)
$data = @(
    @(11, "InternetService_Fiber optic", 2.145842450816183),
    @(28, "PaymentMethod_Electronic check", 1.633358726769917),
    @(6,  "PaperlessBilling", 1.555913962006821),
    @(24, "StreamingMovies_Yes", 1.337121257707364),
    @(22, "StreamingTV_Yes", 1.325961273725601),
    @(1,  "SeniorCitizen", 1.283789169942248),
    @(8,  "TotalCharges", 1.245693478498101),
    @(29, "PaymentMethod_Mailed check", 1.104670690856989),
    @(7,  "MonthlyCharges", 1.088586988520192),
    @(10, "MultipleLines_Yes", 1.042598526618528),
    @(9,  "MultipleLines_No phone service", 1.035111606240906),
    @(2,  "Partner", 1.015554546092614),
    @(27, "PaymentMethod_Credit card (automatic)", 0.9621667868636813),
    @(3,  "Dependents", 0.9473646520053585),
    @(18, "DeviceProtection_Yes", 0.9355378464801896),
    @(0,  "gender", 0.920870414220499),
    @(17, "DeviceProtection_No internet service", 0.8772197782008782),
    @(19, "TechSupport_No internet service", 0.8772197782008782),
    @(21, "StreamingTV_No internet service", 0.8772197782008782),
    @(13, "OnlineSecurity_No internet service", 0.8772197782008782),
    @(23, "StreamingMovies_No internet service", 0.8772197782008782),
    @(12, "InternetService_No", 0.8772197782008782),
    @(15, "OnlineBackup_No internet service", 0.8772197782008782),
    @(16, "OnlineBackup_Yes", 0.7836023713334871),
    @(5,  "PhoneService", 0.756342527765082),
    @(14, "OnlineSecurity_Yes", 0.6952354516484701),
    @(20, "TechSupport_Yes", 0.6593809146412619),
    @(25, "Contract_One year", 0.3660260335714712),
    @(4,  "tenure", 0.3602155390277323),
    @(26, "Contract_Two year", 0.1649222419427319)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
